# Update registration test case documentation: fill in "Actual Outcome" and
# "Fail/Pass" columns now that testing of the registration page is complete.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-10 hold the individual test cases. Column F = Actual Outcome,
# column G = Fail/Pass. Every test case passed with the same actual
# outcome as the expected outcome.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = "Same as expected outcome."
    $ws.Cells.Item($r, 7).Value = "Pass"
}

# Reflect the reviewer's final view state: zoomed to 70%, with the
# selection parked on E10 (the last row reviewed).
[void]$ws.Range("E10").Select()
$excel.ActiveWindow.Zoom = 70
